$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 89.5

# Row 116: Growing Up
$ws.Range("H116").Value = 2847.4
$ws.Range("J116").Value = 4213
$ws.Range("L116").Value = 4213
$ws.Range("N116").Value = -11097

# Row 125: Body over Mind
$ws.Range("H125").Value = 4181.727
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 4444.3335
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 39999.0015
$ws.Range("M125").Value = -24540
$ws.Range("N125").Value = -44919.0015

# Row 129: Practical Command
$ws.Range("H129").Value = 950.4167
$ws.Range("J129").Value = 987.3077
$ws.Range("L129").Value = 2961.9231
$ws.Range("N129").Value = -12961.9231

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2725.8372
$ws.Range("I138").Value = 1357.0435
$ws.Range("J138").Value = 4299.95
$ws.Range("K138").Value = 4071.1305
$ws.Range("L138").Value = 12899.85
$ws.Range("M138").Value = 1068.8695
$ws.Range("N138").Value = -23179.85

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 3640.366
$ws.Range("I74").Value = 4019.4375
$ws.Range("J74").Value = 2292.5557
$ws.Range("K74").Value = 4019.4375
$ws.Range("L74").Value = 2292.5557
$ws.Range("M74").Value = -3145.4375
$ws.Range("N74").Value = -4040.5557

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 3640.366
$ws.Range("I77").Value = 4019.4375
$ws.Range("J77").Value = 2292.5557
$ws.Range("K77").Value = 20097.1875
$ws.Range("L77").Value = 11462.7785
$ws.Range("M77").Value = -15729.1875
$ws.Range("N77").Value = -20198.7785

# Row 103: Sweeping the Legs
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 739.6667
$ws.Range("I122").Value = 684.8
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 2054.4
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = 395.6000000000004
$ws.Range("N122").Value = -7942

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3057.7778
$ws.Range("I132").Value = 2502
$ws.Range("J132").Value = 5003
$ws.Range("K132").Value = 7506
$ws.Range("L132").Value = 15009
$ws.Range("M132").Value = -4976
$ws.Range("N132").Value = -20069

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 3256.2222
$ws.Range("I86").Value = 2461.2
$ws.Range("J86").Value = 4250
$ws.Range("K86").Value = 2461.2
$ws.Range("L86").Value = 4250
$ws.Range("M86").Value = -1338.2
$ws.Range("N86").Value = -6496

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 3256.2222
$ws.Range("I89").Value = 2461.2
$ws.Range("J89").Value = 4250
$ws.Range("K89").Value = 12306
$ws.Range("L89").Value = 21250
$ws.Range("M89").Value = -6690
$ws.Range("N89").Value = -32482

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 2985.8823
$ws.Range("I86").Value = 3632
$ws.Range("J86").Value = 1801.3334
$ws.Range("K86").Value = 3632
$ws.Range("L86").Value = 1801.3334
$ws.Range("M86").Value = -2509
$ws.Range("N86").Value = -4047.3334

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 2985.8823
$ws.Range("I89").Value = 3632
$ws.Range("J89").Value = 1801.3334
$ws.Range("K89").Value = 18160
$ws.Range("L89").Value = 9006.666999999999
$ws.Range("M89").Value = -12544
$ws.Range("N89").Value = -20238.667

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1625
$ws.Range("I122").Value = 1571.4286
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4714.2858
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2264.2858
$ws.Range("N122").Value = -10900

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2690.84
$ws.Range("I134").Value = 1501.2
$ws.Range("K134").Value = 4503.6
$ws.Range("M134").Value = -1968.6

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 996.2432
$ws.Range("J68").Value = 1366.5555
$ws.Range("L68").Value = 4099.666499999999
$ws.Range("N68").Value = -5721.666499999999

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 996.2432
$ws.Range("J71").Value = 1366.5555
$ws.Range("L71").Value = 12298.9995
$ws.Range("N71").Value = -20410.9995

# Row 76: Old Victories, New Tastes
$ws.Range("H76").Value = 100005480
$ws.Range("J76").Value = 6476.875
$ws.Range("L76").Value = 19430.625
$ws.Range("N76").Value = -20196.625

# Row 79: The Eats of Authenticity (L)
$ws.Range("H79").Value = 100005480
$ws.Range("J79").Value = 6476.875
$ws.Range("L79").Value = 19430.625
$ws.Range("N79").Value = -22082.625

# Row 107: Slippery Service
$ws.Range("H107").Value = 38462150
$ws.Range("J107").Value = 100001280
$ws.Range("L107").Value = 300003840
$ws.Range("N107").Value = -300007680

# Row 125: At Any Temperature
$ws.Range("H125").Value = 4691.6665
$ws.Range("I125").Value = 850
$ws.Range("J125").Value = 5460
$ws.Range("K125").Value = 2550
$ws.Range("L125").Value = 16380
$ws.Range("M125").Value = 2370
$ws.Range("N125").Value = -26220

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me
$ws.Range("H5").Value = 10401
$ws.Range("J5").Value = 11668.333
$ws.Range("L5").Value = 11668.333
$ws.Range("N5").Value = -11892.333

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2954.5264
$ws.Range("I122").Value = 2881.5454
$ws.Range("J122").Value = 3054.875
$ws.Range("K122").Value = 8644.636200000001
$ws.Range("L122").Value = 9164.625
$ws.Range("M122").Value = -6194.636200000001
$ws.Range("N122").Value = -14064.625

# Row 132: On Board for Lar
$ws.Range("H132").Value = 4987.6665
$ws.Range("I132").Value = 5073.7354
$ws.Range("J132").Value = 4402.4
$ws.Range("K132").Value = 15221.2062
$ws.Range("L132").Value = 13207.2
$ws.Range("M132").Value = -12691.2062
$ws.Range("N132").Value = -18267.2

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 6800
$ws.Range("J2").Value = 6800
$ws.Range("L2").Value = 6800
$ws.Range("N2").Value = -7024

# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 4654.727
$ws.Range("I7").Value = 4355.8887
$ws.Range("J7").Value = 5999.5
$ws.Range("K7").Value = 4355.8887
$ws.Range("L7").Value = 5999.5
$ws.Range("M7").Value = -4243.8887
$ws.Range("N7").Value = -6223.5

# Row 126: Battered Books
$ws.Range("H126").Value = 4654.727
$ws.Range("I126").Value = 4355.8887
$ws.Range("J126").Value = 5999.5
$ws.Range("K126").Value = 13067.6661
$ws.Range("L126").Value = 17998.5
$ws.Range("M126").Value = -10597.6661
$ws.Range("N126").Value = -22938.5

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4539.2856
$ws.Range("I132").Value = 4986.846
$ws.Range("J132").Value = 3812
$ws.Range("K132").Value = 14960.538
$ws.Range("L132").Value = 11436
$ws.Range("M132").Value = -12430.538
$ws.Range("N132").Value = -16496

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Range("H2").Value = 25326.666
$ws.Range("J2").Value = 25326.666
$ws.Range("L2").Value = 25326.666
$ws.Range("N2").Value = -25550.666
